$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: BigIntegerInstantiation (no rule id / obs, counts 0/0)
$ws.Range("B6").Value = "category/java/performance.xml/BigIntegerInstantiation"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

# Row 7: BooleanInstantiation
$ws.Range("A7").Value = "BooleanInstantiation"
$ws.Range("B7").Value = "category/java/performance.xml/BooleanInstantiation"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1

# Row 8: ByteInstantiation
$ws.Range("B8").Value = "category/java/performance.xml/ByteInstantiation"
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1

# Select entire sheet (Ctrl+A twice), mirroring the end-user action that
# produced sqref="A1:XFD1048576" in the saved view state.
$ws.Cells.Select()
